$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from AE1 (DiffA) into the new header cells AF1:AH1
$ws.Range("AE1").Copy($ws.Range("AF1"))
$ws.Range("AE1").Copy($ws.Range("AG1"))
$ws.Range("AE1").Copy($ws.Range("AH1"))

$ws.Range("AF1").Value = "%DiffH"
$ws.Range("AG1").Value = "%DiffD"
$ws.Range("AH1").Value = "%DiffA"

# Populate %DiffH, %DiffD, %DiffA = DiffX / YtrueX * 100 for each data row
$ws.Range("AF2").Value = 261.8808734230747
$ws.Range("AG2").Value = -90.36177304850587
$ws.Range("AH2").Value = -34.98396227006597
$ws.Range("AF3").Value = 193.5067766979248
$ws.Range("AG3").Value = -91.39852989904118
$ws.Range("AH3").Value = -84.34066324358595
$ws.Range("AF4").Value = 189.96237896973426
$ws.Range("AG4").Value = -93.2017343919541
$ws.Range("AH4").Value = -53.818223937403296
$ws.Range("AF5").Value = 376.146181719355
$ws.Range("AG5").Value = -89.64708157335758
$ws.Range("AH5").Value = -73.16398145128963
$ws.Range("AF6").Value = 199.23201757512837
$ws.Range("AG6").Value = -88.32753755093583
$ws.Range("AH6").Value = -77.14487091161547
$ws.Range("AF7").Value = 256.8211614598018
$ws.Range("AG7").Value = -86.3666869749245
$ws.Range("AH7").Value = -53.784816824646974
$ws.Range("AF8").Value = 228.3609094170749
$ws.Range("AG8").Value = -92.88601650191526
$ws.Range("AH8").Value = -74.07905406016857
$ws.Range("AF9").Value = 211.6354706897885
$ws.Range("AG9").Value = -89.74475817717703
$ws.Range("AH9").Value = -74.72340595787425
$ws.Range("AF10").Value = 339.0720019100143
$ws.Range("AG10").Value = -90.9752779679102
$ws.Range("AH10").Value = -51.71372408027244
$ws.Range("AF11").Value = 224.14995811707223
$ws.Range("AG11").Value = -99.43319744862754
$ws.Range("AH11").Value = -51.68859345704882
$ws.Range("AF12").Value = 311.705692251894
$ws.Range("AG12").Value = -93.8398723215762
$ws.Range("AH12").Value = -52.307969750459236
$ws.Range("AF13").Value = 277.0747178228872
$ws.Range("AG13").Value = -94.62539114824929
$ws.Range("AH13").Value = -46.46970822871661
$ws.Range("AF14").Value = 234.02137130202743
$ws.Range("AG14").Value = -88.75726544523889
$ws.Range("AH14").Value = -44.08643482539234
$ws.Range("AF15").Value = 290.72783722523695
$ws.Range("AG15").Value = -87.77403788631516
$ws.Range("AH15").Value = -51.08344690446458
$ws.Range("AF16").Value = 252.71452304949068
$ws.Range("AG16").Value = -93.06265820652769
$ws.Range("AH16").Value = -29.787885419384402
$ws.Range("AF17").Value = 271.96738632530423
$ws.Range("AG17").Value = -90.51522160416094
$ws.Range("AH17").Value = -43.934900944481
$ws.Range("AF18").Value = 131.53244072187965
$ws.Range("AG18").Value = -96.83681357671509
$ws.Range("AH18").Value = -71.12339073805352
$ws.Range("AF19").Value = 276.20074062114656
$ws.Range("AG19").Value = -91.71859035245284
$ws.Range("AH19").Value = -27.78480940497407
$ws.Range("AF20").Value = 180.99785712612487
$ws.Range("AG20").Value = -95.27854842299925
$ws.Range("AH20").Value = -39.31099669572695
$ws.Range("AF21").Value = 261.8897941494533
$ws.Range("AG21").Value = -82.66340228101768
$ws.Range("AH21").Value = -68.35626731399724
$ws.Range("AF22").Value = 214.18210447573304
$ws.Range("AG22").Value = -90.37941950060731
$ws.Range("AH22").Value = -46.50870250355024
$ws.Range("AF23").Value = 176.59825493683766
$ws.Range("AG23").Value = -90.84717698752024
$ws.Range("AH23").Value = -76.26268766776583
$ws.Range("AF24").Value = 242.24261430097482
$ws.Range("AG24").Value = -88.44235674390052
$ws.Range("AH24").Value = -55.91543904870636
$ws.Range("AF25").Value = 141.236184116541
$ws.Range("AG25").Value = -92.42999412255605
$ws.Range("AH25").Value = -50.982588090443535
$ws.Range("AF26").Value = 275.8830799435723
$ws.Range("AG26").Value = -91.61673194677057
$ws.Range("AH26").Value = -34.90382788273811
$ws.Range("AF27").Value = 367.40198089291283
$ws.Range("AG27").Value = -89.2463413656163
$ws.Range("AH27").Value = -59.31013790089549
$ws.Range("AF28").Value = 210.81827379374906
$ws.Range("AG28").Value = -89.2421340647114
$ws.Range("AH28").Value = -31.23297279904978
$ws.Range("AF29").Value = 242.63501602464825
$ws.Range("AG29").Value = -89.41497700893134
$ws.Range("AH29").Value = -41.9666850954375
$ws.Range("AF30").Value = 322.12880595418864
$ws.Range("AG30").Value = -87.57071627901813
$ws.Range("AH30").Value = -44.7545117511434
$ws.Range("AF31").Value = 329.11557699415386
$ws.Range("AG31").Value = -85.92549778231641
$ws.Range("AH31").Value = -43.97530844484992
$ws.Range("AF32").Value = 194.50528747764787
$ws.Range("AG32").Value = -95.29776635192152
$ws.Range("AH32").Value = -39.06261590513966
$ws.Range("AF33").Value = 196.8572419465419
$ws.Range("AG33").Value = -88.25182925843784
$ws.Range("AH33").Value = -39.81373284404978
$ws.Range("AF34").Value = 237.2841103314693
$ws.Range("AG34").Value = -87.39030151799602
$ws.Range("AH34").Value = -45.791902082378655
$ws.Range("AF35").Value = 138.7212513060598
$ws.Range("AG35").Value = -86.25739678093531
$ws.Range("AH35").Value = -72.22711403888289
$ws.Range("AF36").Value = 270.25540038586433
$ws.Range("AG36").Value = -92.06770804293187
$ws.Range("AH36").Value = -39.74470862931711
$ws.Range("AF37").Value = 321.3234363524414
$ws.Range("AG37").Value = -90.30526312158544
$ws.Range("AH37").Value = -57.96141144751466
$ws.Range("AF38").Value = 325.35581520139357
$ws.Range("AG38").Value = -89.86913852753439
$ws.Range("AH38").Value = -41.65664756921406
$ws.Range("AF39").Value = 272.53472132293905
$ws.Range("AG39").Value = -94.75176009303715
$ws.Range("AH39").Value = -30.7420311204459
$ws.Range("AF40").Value = 250.4258121578041
$ws.Range("AG40").Value = -94.80429110104143
$ws.Range("AH40").Value = -43.973327410693784
$ws.Range("AF41").Value = 154.38825370788
$ws.Range("AG41").Value = -89.25646941425441
$ws.Range("AH41").Value = -64.07577187911359
$ws.Range("AF42").Value = 132.3123128041287
$ws.Range("AG42").Value = -86.93828778894166
$ws.Range("AH42").Value = -82.81855035184391
$ws.Range("AF43").Value = 302.5597652971716
$ws.Range("AG43").Value = -90.31276203718144
$ws.Range("AH43").Value = -40.70025732320578
